$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date style (numFmt 14, font 1) from A172 onto the new date cells A173:A235
$ws.Range("A172").Copy()
$ws.Range("A173:A235").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A173").Value = 44022
$ws.Range("B173").Value = 0.691150477323836
$ws.Range("A174").Value = 44023
$ws.Range("B174").Value = 0.749834665457363
$ws.Range("A175").Value = 44024
$ws.Range("B175").Value = 0.740101918541116
$ws.Range("A176").Value = 44025
$ws.Range("B176").Value = 0.673262004386364
$ws.Range("A177").Value = 44026
$ws.Range("B177").Value = 0.709073935062636
$ws.Range("A178").Value = 44027
$ws.Range("B178").Value = 0.672091426247228
$ws.Range("A179").Value = 44028
$ws.Range("B179").Value = 0.703076697580215
$ws.Range("A180").Value = 44029
$ws.Range("B180").Value = 0.662027573911487
$ws.Range("A181").Value = 44030
$ws.Range("B181").Value = 0.552484794470793
$ws.Range("A182").Value = 44031
$ws.Range("B182").Value = 0.679105939963991
$ws.Range("A183").Value = 44032
$ws.Range("B183").Value = 0.740690333527513
$ws.Range("A184").Value = 44033
$ws.Range("B184").Value = 0.695784054897994
$ws.Range("A185").Value = 44034
$ws.Range("B185").Value = 0.709039548022599
$ws.Range("A186").Value = 44035
$ws.Range("B186").Value = 0.703610343900821
$ws.Range("A187").Value = 44036
$ws.Range("B187").Value = 0.703663958813894
$ws.Range("A188").Value = 44037
$ws.Range("B188").Value = 0.730821158319912
$ws.Range("A189").Value = 44038
$ws.Range("B189").Value = 0.756914913060961
$ws.Range("A190").Value = 44039
$ws.Range("B190").Value = 0.676236923022429
$ws.Range("A191").Value = 44040
$ws.Range("B191").Value = 0.61715158169217
$ws.Range("A192").Value = 44041
$ws.Range("B192").Value = 0.605839500887566
$ws.Range("A193").Value = 44042
$ws.Range("B193").Value = 0.710039654589625
$ws.Range("A194").Value = 44043
$ws.Range("B194").Value = 0.707875408128016
$ws.Range("A195").Value = 44044
$ws.Range("B195").Value = 0.706420864601834
$ws.Range("A196").Value = 44045
$ws.Range("B196").Value = 0.704788678835289
$ws.Range("A197").Value = 44046
$ws.Range("B197").Value = 0.665881890816702
$ws.Range("A198").Value = 44047
$ws.Range("B198").Value = 0.716048699740972
$ws.Range("A199").Value = 44048
$ws.Range("B199").Value = 0.628678126885428
$ws.Range("A200").Value = 44049
$ws.Range("B200").Value = 0.725564292698827
$ws.Range("A201").Value = 44050
$ws.Range("B201").Value = 0.72798697923324
$ws.Range("A202").Value = 44051
$ws.Range("B202").Value = 0.659954015971945
$ws.Range("A203").Value = 44052
$ws.Range("B203").Value = 0.560796280477753
$ws.Range("A204").Value = 44053
$ws.Range("B204").Value = 0.682123570508506
$ws.Range("A205").Value = 44054
$ws.Range("B205").Value = 0.62735539030266
$ws.Range("A206").Value = 44055
$ws.Range("B206").Value = 0.660660171253392
$ws.Range("A207").Value = 44056
$ws.Range("B207").Value = 0.68769227910062
$ws.Range("A208").Value = 44057
$ws.Range("B208").Value = 0.635008401621019
$ws.Range("A209").Value = 44058
$ws.Range("B209").Value = 0.616470428852726
$ws.Range("A210").Value = 44059
$ws.Range("B210").Value = 0.740133365927857
$ws.Range("A211").Value = 44060
$ws.Range("B211").Value = 0.696402335590417
$ws.Range("A212").Value = 44061
$ws.Range("B212").Value = 0.740238706136685
$ws.Range("A213").Value = 44062
$ws.Range("B213").Value = 0.67089678863029
$ws.Range("A214").Value = 44063
$ws.Range("B214").Value = 0.723626420660319
$ws.Range("A215").Value = 44064
$ws.Range("B215").Value = 0.737865821520708
$ws.Range("A216").Value = 44065
$ws.Range("B216").Value = 0.617668207498716
$ws.Range("A217").Value = 44066
$ws.Range("B217").Value = 0.703520502249316
$ws.Range("A218").Value = 44067
$ws.Range("B218").Value = 0.731665469316801
$ws.Range("A219").Value = 44068
$ws.Range("B219").Value = 0.679070473071093
$ws.Range("A220").Value = 44069
$ws.Range("B220").Value = 0.648828884466708
$ws.Range("A221").Value = 44070
$ws.Range("B221").Value = 0.574554069119287
$ws.Range("A222").Value = 44071
$ws.Range("B222").Value = 0.784280322228475
$ws.Range("A223").Value = 44072
$ws.Range("B223").Value = 0.929595827900913
$ws.Range("A224").Value = 44073
$ws.Range("B224").Value = 0.704621949753576
$ws.Range("A225").Value = 44074
$ws.Range("B225").Value = 0.680320338307517
$ws.Range("A226").Value = 44075
$ws.Range("B226").Value = 0.721732027027758
$ws.Range("A227").Value = 44076
$ws.Range("B227").Value = 0.71068386415844
$ws.Range("A228").Value = 44077
$ws.Range("B228").Value = 0.676447206571677
$ws.Range("A229").Value = 44078
$ws.Range("B229").Value = 0.726232709295262
$ws.Range("A230").Value = 44079
$ws.Range("B230").Value = 0.727504097394426
$ws.Range("A231").Value = 44080
$ws.Range("B231").Value = 0.625196081456737
$ws.Range("A232").Value = 44081
$ws.Range("B232").Value = 0.748234402647141
$ws.Range("A233").Value = 44082
$ws.Range("B233").Value = 0.714552772866274
$ws.Range("A234").Value = 44083
$ws.Range("B234").Value = 0.722560385710641
$ws.Range("A235").Value = 44084
$ws.Range("B235").Value = 0.718145430390107
